# Move "delivery_organisation_path" (currently the last data column, R) on the
# "Service Contacts" sheet to sit right before "practitioner_key" (column D),
# i.e. immediately after "episode_key" (column C).
#
# This is equivalent to the user selecting column R, cutting it, and then using
# "Insert Cut Cells" on column D - which shifts the old D:Q block one column to
# the right (into E:R) and drops the cut column into the vacated D.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Service Contacts")
$wsK5 = $wb.Worksheets.Item("K5")

$ws.Columns("R:R").Cut()
$ws.Columns("D:D").Insert()

# The insert-shift can leave behind stray "not actually custom" width markers on
# the columns that moved from D:Q to E:R; clear them so only the genuinely
# custom-width columns (the relocated D, plus Q:R which inherited the old P:Q
# widths) keep explicit formatting.
$ws.Columns("E:Q").ClearFormats()

# Reflect the resulting selection state: the newly inserted column D ends up
# selected on Service Contacts, and the user's last action was selecting column
# F on the K5 sheet (which stays the active tab).
$ws.Range("D1:D1048576").Select()
$wsK5.Activate()
$wsK5.Range("F1:F5").Select()
